# Pooh Points: normal 20260214
# - Narrow the "status" column (G) on the "Players" sheet from 17 to 8
# - Update every in-progress game-clock status (e.g. "1:51 - 2nd Half",
#   "0:13 - 2nd Half") in column G to "Final"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Narrow column G (status) so the stored OOXML width becomes 8.
# Excel's ColumnWidth character units are offset from the stored
# worksheet width by 5/6, so subtract that to land exactly on 8.
$ws.Columns.Item(7).ColumnWidth = 8 - 5/6

# Walk every used row and flip any non-final game clock to "Final"
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value
    if ($val -ne $null -and $val -ne "Final") {
        $cell.Value = "Final"
    }
}
